# Auto-generated script to apply market-data refresh values
# as described by the commit "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 956.5
$ws.Range("I19").Value = 916.6667
$ws.Range("K19").Value = 916.6667
$ws.Range("M19").Value = -741.6667

$ws.Range("H33").Value = 119.181816
$ws.Range("I33").Value = 127.888885
$ws.Range("K33").Value = 127.888885
$ws.Range("M33").Value = 101.111115

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H76").Value = 3535
$ws.Range("I76").Value = 2950
$ws.Range("J76").Value = 3925
$ws.Range("K76").Value = 2950
$ws.Range("L76").Value = 3925
$ws.Range("M76").Value = -2635
$ws.Range("N76").Value = -4555

$ws.Range("H79").Value = 3535
$ws.Range("I79").Value = 2950
$ws.Range("J79").Value = 3925
$ws.Range("K79").Value = 2950
$ws.Range("L79").Value = 3925
$ws.Range("M79").Value = -1858
$ws.Range("N79").Value = -6109

$ws.Range("H92").Value = 597
$ws.Range("I92").Value = 643.8889
$ws.Range("J92").Value = 476.42856
$ws.Range("K92").Value = 643.8889
$ws.Range("L92").Value = 476.42856
$ws.Range("M92").Value = 604.1111
$ws.Range("N92").Value = -2972.42856

$ws.Range("H96").Value = 513.44446
$ws.Range("I96").Value = 327.625
$ws.Range("K96").Value = 982.875
$ws.Range("M96").Value = 390.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 13163.917
$ws.Range("I74").Value = 12335.8
$ws.Range("J74").Value = 17304.5
$ws.Range("K74").Value = 12335.8
$ws.Range("L74").Value = 17304.5
$ws.Range("M74").Value = -11461.8
$ws.Range("N74").Value = -19052.5

$ws.Range("H77").Value = 13163.917
$ws.Range("I77").Value = 12335.8
$ws.Range("J77").Value = 17304.5
$ws.Range("K77").Value = 61679
$ws.Range("L77").Value = 86522.5
$ws.Range("M77").Value = -57311
$ws.Range("N77").Value = -95258.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3097
$ws.Range("I94").Value = 3620
$ws.Range("J94").Value = 743.5
$ws.Range("K94").Value = 3620
$ws.Range("L94").Value = 743.5
$ws.Range("M94").Value = -3169
$ws.Range("N94").Value = -1645.5

$ws.Range("H134").Value = 2935.0715
$ws.Range("I134").Value = 2776.2307
$ws.Range("K134").Value = 8328.6921
$ws.Range("M134").Value = -5793.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

$ws.Range("H62").Value = 8999.799999999999
$ws.Range("J62").Value = 8749.75
$ws.Range("L62").Value = 8749.75
$ws.Range("N62").Value = -9997.75

$ws.Range("H65").Value = 8999.799999999999
$ws.Range("J65").Value = 8749.75
$ws.Range("L65").Value = 43748.75
$ws.Range("N65").Value = -49988.75

$ws.Range("H86").Value = 3589.7222
$ws.Range("I86").Value = 4008.4443
$ws.Range("K86").Value = 4008.4443
$ws.Range("M86").Value = -2885.4443

$ws.Range("H89").Value = 3589.7222
$ws.Range("I89").Value = 4008.4443
$ws.Range("K89").Value = 20042.2215
$ws.Range("M89").Value = -14426.2215

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1003.6667
$ws.Range("I64").Value = 1003.6667
$ws.Range("K64").Value = 3011.0001
$ws.Range("M64").Value = -2741.0001

$ws.Range("H67").Value = 1003.6667
$ws.Range("I67").Value = 1003.6667
$ws.Range("K67").Value = 3011.0001
$ws.Range("M67").Value = -2075.0001

$ws.Range("H132").Value = 3122.375
$ws.Range("I132").Value = 1195.8
$ws.Range("K132").Value = 10762.2
$ws.Range("M132").Value = -8232.199999999999

$ws.Range("H137").Value = 2786.3333
$ws.Range("I137").Value = 2274.75
$ws.Range("K137").Value = 6824.25
$ws.Range("M137").Value = -1724.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2928.4443
$ws.Range("I80").Value = 2748.5
$ws.Range("J80").Value = 3072.4
$ws.Range("K80").Value = 2748.5
$ws.Range("L80").Value = 3072.4
$ws.Range("M80").Value = -1750.5
$ws.Range("N80").Value = -5068.4

$ws.Range("H83").Value = 2928.4443
$ws.Range("I83").Value = 2748.5
$ws.Range("J83").Value = 3072.4
$ws.Range("K83").Value = 13742.5
$ws.Range("L83").Value = 15362
$ws.Range("M83").Value = -8750.5
$ws.Range("N83").Value = -25346

$ws.Range("H97").Value = 2011
$ws.Range("I97").Value = 2155.4285
$ws.Range("K97").Value = 2155.4285
$ws.Range("M97").Value = -1659.4285

$ws.Range("H132").Value = 2895.8635
$ws.Range("I132").Value = 2247.5386
$ws.Range("K132").Value = 6742.6158
$ws.Range("M132").Value = -4212.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2889.9285
$ws.Range("I68").Value = 2889.9285
$ws.Range("K68").Value = 2889.9285
$ws.Range("M68").Value = -2140.9285

$ws.Range("H71").Value = 2889.9285
$ws.Range("I71").Value = 2889.9285
$ws.Range("K71").Value = 14449.6425
$ws.Range("M71").Value = -10705.6425

$ws.Range("H132").Value = 5907.909
$ws.Range("J132").Value = 5907.909
$ws.Range("L132").Value = 17723.727
$ws.Range("N132").Value = -22783.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4769.6
$ws.Range("I62").Value = 3449.75
$ws.Range("K62").Value = 3449.75
$ws.Range("M62").Value = -2825.75

$ws.Range("H65").Value = 4769.6
$ws.Range("I65").Value = 3449.75
$ws.Range("K65").Value = 17248.75
$ws.Range("M65").Value = -14128.75

$ws.Range("H81").Value = 1527.1111
$ws.Range("I81").Value = 1677.7142
$ws.Range("K81").Value = 3355.4284
$ws.Range("M81").Value = -2294.4284

$ws.Range("H84").Value = 1527.1111
$ws.Range("I84").Value = 1677.7142
$ws.Range("K84").Value = 16777.142
$ws.Range("M84").Value = -11473.142

$ws.Range("H132").Value = 3110.875
$ws.Range("J132").Value = 3099.182
$ws.Range("L132").Value = 9297.545999999998
$ws.Range("N132").Value = -14357.546

$ws.Range("H136").Value = 4920.0835
$ws.Range("I136").Value = 4255.5
$ws.Range("J136").Value = 6249.25
$ws.Range("K136").Value = 12766.5
$ws.Range("L136").Value = 18747.75
$ws.Range("M136").Value = -10216.5
$ws.Range("N136").Value = -23847.75

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
